$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update J2:J34 phone numbers (also clear the per-cell border/center style to plain Normal) ---
$ws.Cells.Item(2, 10).Value = 5616677351
$ws.Cells.Item(2, 10).Style = "Normal"
$ws.Cells.Item(3, 10).Value = 5614858433
$ws.Cells.Item(3, 10).Style = "Normal"
$ws.Cells.Item(4, 10).Value = 5615869888
$ws.Cells.Item(4, 10).Style = "Normal"
$ws.Cells.Item(5, 10).Value = 5614714457
$ws.Cells.Item(5, 10).Style = "Normal"
$ws.Cells.Item(6, 10).Value = 5619647774
$ws.Cells.Item(6, 10).Style = "Normal"
$ws.Cells.Item(7, 10).Value = 5614786795
$ws.Cells.Item(7, 10).Style = "Normal"
$ws.Cells.Item(8, 10).Value = 5619744199
$ws.Cells.Item(8, 10).Style = "Normal"
$ws.Cells.Item(9, 10).Value = 5618645137
$ws.Cells.Item(9, 10).Style = "Normal"
$ws.Cells.Item(10, 10).Value = 5616914758
$ws.Cells.Item(10, 10).Style = "Normal"
$ws.Cells.Item(11, 10).Value = 5616667133
$ws.Cells.Item(11, 10).Style = "Normal"
$ws.Cells.Item(12, 10).Value = 5615945674
$ws.Cells.Item(12, 10).Style = "Normal"
$ws.Cells.Item(13, 10).Value = 5615853684
$ws.Cells.Item(13, 10).Style = "Normal"
$ws.Cells.Item(14, 10).Value = 5615938591
$ws.Cells.Item(14, 10).Style = "Normal"
$ws.Cells.Item(15, 10).Value = 5616677373
$ws.Cells.Item(15, 10).Style = "Normal"
$ws.Cells.Item(16, 10).Value = 5619831869
$ws.Cells.Item(16, 10).Style = "Normal"
$ws.Cells.Item(17, 10).Value = 5618956315
$ws.Cells.Item(17, 10).Style = "Normal"
$ws.Cells.Item(18, 10).Value = 5616914654
$ws.Cells.Item(18, 10).Style = "Normal"
$ws.Cells.Item(19, 10).Value = 5618954116
$ws.Cells.Item(19, 10).Style = "Normal"
$ws.Cells.Item(20, 10).Value = 5618956315
$ws.Cells.Item(20, 10).Style = "Normal"
$ws.Cells.Item(21, 10).Value = 5613917133
$ws.Cells.Item(21, 10).Style = "Normal"
$ws.Cells.Item(22, 10).Value = 5616695198
$ws.Cells.Item(22, 10).Style = "Normal"
$ws.Cells.Item(23, 10).Value = 5615955997
$ws.Cells.Item(23, 10).Style = "Normal"
$ws.Cells.Item(24, 10).Value = 5614854416
$ws.Cells.Item(24, 10).Style = "Normal"
$ws.Cells.Item(25, 10).Value = 5613643175
$ws.Cells.Item(25, 10).Style = "Normal"
$ws.Cells.Item(26, 10).Value = 5614788353
$ws.Cells.Item(26, 10).Style = "Normal"
$ws.Cells.Item(27, 10).Value = 5616921415
$ws.Cells.Item(27, 10).Style = "Normal"
$ws.Cells.Item(28, 10).Value = 5615749433
$ws.Cells.Item(28, 10).Style = "Normal"
$ws.Cells.Item(29, 10).Value = 5619843437
$ws.Cells.Item(29, 10).Style = "Normal"
$ws.Cells.Item(30, 10).Value = 5619813471
$ws.Cells.Item(30, 10).Style = "Normal"
$ws.Cells.Item(31, 10).Value = 5614975479
$ws.Cells.Item(31, 10).Style = "Normal"
$ws.Cells.Item(32, 10).Value = 5617755983
$ws.Cells.Item(32, 10).Style = "Normal"
$ws.Cells.Item(33, 10).Value = 5616744552
$ws.Cells.Item(33, 10).Style = "Normal"
$ws.Cells.Item(34, 10).Value = 5617671652
$ws.Cells.Item(34, 10).Style = "Normal"

# --- Append 21 new rows (97-117) for Alborz / Karaj district 4 ---
# Copy the bordered/centered formatting used by existing data rows onto the new block first,
# so the new cells reuse the same cell style (no new style entries are introduced).
$ws.Range("A2:K2").Copy()
$ws.Range("A97:K117").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 97
$ws.Cells.Item(97, 1).Value = "البرز"
$ws.Cells.Item(97, 2).Value = "کرج ناحيه 4"
$ws.Cells.Item(97, 3).Value = "شهيدان  هداوند"
$ws.Cells.Item(97, 4).Value = "دوره ابتدايي توصيفي"
$ws.Cells.Item(97, 5).Value = "دخترانه"
$ws.Cells.Item(97, 6).Value = "دولتي"
$ws.Cells.Item(97, 7).Value = "عادي"
$ws.Cells.Item(97, 8).Value = "عادي"
$ws.Cells.Item(97, 9).Value = 4552018
$ws.Cells.Item(97, 10).Value = 319777333
$ws.Cells.Item(97, 11).Value = "حصار ک بالا خيابان آقارضايي آموزشگاه شهيدان هداوند2"

# Row 98
$ws.Cells.Item(98, 1).Value = "البرز"
$ws.Cells.Item(98, 2).Value = "کرج ناحيه 4"
$ws.Cells.Item(98, 3).Value = "شهداي کمالشهر"
$ws.Cells.Item(98, 4).Value = "دوره متوسطه اول"
$ws.Cells.Item(98, 5).Value = "پسرانه"
$ws.Cells.Item(98, 6).Value = "دولتي"
$ws.Cells.Item(98, 7).Value = "عادي"
$ws.Cells.Item(98, 8).Value = "عادي"
$ws.Cells.Item(98, 9).Value = 4703847
$ws.Cells.Item(98, 10).Value = 319976577
$ws.Cells.Item(98, 11).Value = "کمالشهر_خ وليعصرجنوبي _خ شهدا"

# Row 99
$ws.Cells.Item(99, 1).Value = "البرز"
$ws.Cells.Item(99, 2).Value = "کرج ناحيه 4"
$ws.Cells.Item(99, 3).Value = "شهيد صياد شيرازي"
$ws.Cells.Item(99, 4).Value = "دوره ابتدايي توصيفي"
$ws.Cells.Item(99, 5).Value = "پسرانه"
$ws.Cells.Item(99, 6).Value = "دولتي"
$ws.Cells.Item(99, 7).Value = "عادي"
$ws.Cells.Item(99, 8).Value = "عادي"
$ws.Cells.Item(99, 9).Value = 34703847
$ws.Cells.Item(99, 10).Value = "-"
$ws.Cells.Item(99, 11).Value = "کمالشهر جنب مسجد علي ابن ابي طالب- مدرسه ش صياد شيرازي"

# Row 100
$ws.Cells.Item(100, 1).Value = "البرز"
$ws.Cells.Item(100, 2).Value = "کرج ناحيه 4"
$ws.Cells.Item(100, 3).Value = "شهيدستارلطفي"
$ws.Cells.Item(100, 4).Value = "متوسطه دوم - نظري"
$ws.Cells.Item(100, 5).Value = "پسرانه"
$ws.Cells.Item(100, 6).Value = "دولتي"
$ws.Cells.Item(100, 7).Value = "عادي"
$ws.Cells.Item(100, 8).Value = "عادي"
$ws.Cells.Item(100, 9).Value = "-"
$ws.Cells.Item(100, 10).Value = "-"
$ws.Cells.Item(100, 11).Value = "کيانمهر"

# Row 101
$ws.Cells.Item(101, 1).Value = "البرز"
$ws.Cells.Item(101, 2).Value = "کرج ناحيه 4"
$ws.Cells.Item(101, 3).Value = "هنرستان امام خميني (ره )(1)"
$ws.Cells.Item(101, 4).Value = "متوسطه دوم - هنرستان فني"
$ws.Cells.Item(101, 5).Value = "پسرانه"
$ws.Cells.Item(101, 6).Value = "دولتي"
$ws.Cells.Item(101, 7).Value = "عادي"
$ws.Cells.Item(101, 8).Value = "هيات امنايي"
$ws.Cells.Item(101, 9).Value = 3528485
$ws.Cells.Item(101, 10).Value = 318383464
$ws.Cells.Item(101, 11).Value = "فاز4مهرشهر_بلوار گلها _چهارراه هنرستان خ 406شرقي-پ192"

# Row 102
$ws.Cells.Item(102, 1).Value = "البرز"
$ws.Cells.Item(102, 2).Value = "کرج ناحيه 4"
$ws.Cells.Item(102, 3).Value = "حضرت مريم (1)"
$ws.Cells.Item(102, 4).Value = "متوسطه دوم - نظري"
$ws.Cells.Item(102, 5).Value = "دخترانه"
$ws.Cells.Item(102, 6).Value = "دولتي"
$ws.Cells.Item(102, 7).Value = "عادي"
$ws.Cells.Item(102, 8).Value = "عادي"
$ws.Cells.Item(102, 9).Value = 3408607
$ws.Cells.Item(102, 10).Value = 318579855
$ws.Cells.Item(102, 11).Value = "مهرشهر_بلوارارم _بلواردانش _خ 100_خ مريم"

# Row 103
$ws.Cells.Item(103, 1).Value = "البرز"
$ws.Cells.Item(103, 2).Value = "کرج ناحيه 4"
$ws.Cells.Item(103, 3).Value = "شهيد شهسواري(2)"
$ws.Cells.Item(103, 4).Value = "دوره متوسطه اول"
$ws.Cells.Item(103, 5).Value = "پسرانه"
$ws.Cells.Item(103, 6).Value = "دولتي"
$ws.Cells.Item(103, 7).Value = "عادي"
$ws.Cells.Item(103, 8).Value = "عادي"
$ws.Cells.Item(103, 9).Value = 34662060
$ws.Cells.Item(103, 10).Value = "-"
$ws.Cells.Item(103, 11).Value = "کرج- حصارک پايين -رضاشهر-انتهاي خيابان فروردين"

# Row 104
$ws.Cells.Item(104, 1).Value = "البرز"
$ws.Cells.Item(104, 2).Value = "کرج ناحيه 4"
$ws.Cells.Item(104, 3).Value = "لقمان حکيم(1)"
$ws.Cells.Item(104, 4).Value = "دوره متوسطه اول"
$ws.Cells.Item(104, 5).Value = "پسرانه"
$ws.Cells.Item(104, 6).Value = "دولتي"
$ws.Cells.Item(104, 7).Value = "عادي"
$ws.Cells.Item(104, 8).Value = "عادي"
$ws.Cells.Item(104, 9).Value = 34801188
$ws.Cells.Item(104, 10).Value = "-"
$ws.Cells.Item(104, 11).Value = "کرج-پيشاهنگي-گلدشت"

# Row 105
$ws.Cells.Item(105, 1).Value = "البرز"
$ws.Cells.Item(105, 2).Value = "کرج ناحيه 4"
$ws.Cells.Item(105, 3).Value = "هدايت"
$ws.Cells.Item(105, 4).Value = "دوره ابتدايي توصيفي"
$ws.Cells.Item(105, 5).Value = "دخترانه"
$ws.Cells.Item(105, 6).Value = "دولتي"
$ws.Cells.Item(105, 7).Value = "عادي"
$ws.Cells.Item(105, 8).Value = "عادي"
$ws.Cells.Item(105, 9).Value = "-"
$ws.Cells.Item(105, 10).Value = "-"
$ws.Cells.Item(105, 11).Value = "-"

# Row 106
$ws.Cells.Item(106, 1).Value = "البرز"
$ws.Cells.Item(106, 2).Value = "کرج ناحيه 4"
$ws.Cells.Item(106, 3).Value = "شهيد پرورش"
$ws.Cells.Item(106, 4).Value = "دوره ابتدايي توصيفي"
$ws.Cells.Item(106, 5).Value = "پسرانه"
$ws.Cells.Item(106, 6).Value = "دولتي"
$ws.Cells.Item(106, 7).Value = "عادي"
$ws.Cells.Item(106, 8).Value = "عادي"
$ws.Cells.Item(106, 9).Value = 33214848
$ws.Cells.Item(106, 10).Value = 123454
$ws.Cells.Item(106, 11).Value = "شهرک کيان مهر-خ نبرد اهواز - مدرسه شهيدمحمدپرورشي"

# Row 107
$ws.Cells.Item(107, 1).Value = "البرز"
$ws.Cells.Item(107, 2).Value = "کرج ناحيه 4"
$ws.Cells.Item(107, 3).Value = "امام رضا(ع )2"
$ws.Cells.Item(107, 4).Value = "دوره ابتدايي توصيفي"
$ws.Cells.Item(107, 5).Value = "دخترانه"
$ws.Cells.Item(107, 6).Value = "دولتي"
$ws.Cells.Item(107, 7).Value = "عادي"
$ws.Cells.Item(107, 8).Value = "عادي"
$ws.Cells.Item(107, 9).Value = 3315356
$ws.Cells.Item(107, 10).Value = 318695759
$ws.Cells.Item(107, 11).Value = "جاده قزلحصارروبروي بي سيم شهرک سهرابيه"

# Row 108
$ws.Cells.Item(108, 1).Value = "البرز"
$ws.Cells.Item(108, 2).Value = "کرج ناحيه 4"
$ws.Cells.Item(108, 3).Value = "شهيد باهنر"
$ws.Cells.Item(108, 4).Value = "دوره ابتدايي توصيفي"
$ws.Cells.Item(108, 5).Value = "دخترانه"
$ws.Cells.Item(108, 6).Value = "دولتي"
$ws.Cells.Item(108, 7).Value = "عادي"
$ws.Cells.Item(108, 8).Value = "عادي"
$ws.Cells.Item(108, 9).Value = 33203180
$ws.Cells.Item(108, 10).Value = 318761748
$ws.Cells.Item(108, 11).Value = "کيان مهر ميدان مهرگان بوستان هفتم"

# Row 109
$ws.Cells.Item(109, 1).Value = "البرز"
$ws.Cells.Item(109, 2).Value = "کرج ناحيه 4"
$ws.Cells.Item(109, 3).Value = "وحدت اسلامي (1)"
$ws.Cells.Item(109, 4).Value = "متوسطه دوم - هنرستان کاردانش"
$ws.Cells.Item(109, 5).Value = "پسرانه"
$ws.Cells.Item(109, 6).Value = "دولتي"
$ws.Cells.Item(109, 7).Value = "عادي"
$ws.Cells.Item(109, 8).Value = "عادي"
$ws.Cells.Item(109, 9).Value = 4553805
$ws.Cells.Item(109, 10).Value = 319767745
$ws.Cells.Item(109, 11).Value = "حصارک بالا_روبروي مجتمع ورزشي ايثار"

# Row 110
$ws.Cells.Item(110, 1).Value = "البرز"
$ws.Cells.Item(110, 2).Value = "کرج ناحيه 4"
$ws.Cells.Item(110, 3).Value = "وليعصر(عج)"
$ws.Cells.Item(110, 4).Value = "دوره متوسطه اول"
$ws.Cells.Item(110, 5).Value = "پسرانه"
$ws.Cells.Item(110, 6).Value = "دولتي"
$ws.Cells.Item(110, 7).Value = "عادي"
$ws.Cells.Item(110, 8).Value = "عادي"
$ws.Cells.Item(110, 9).Value = 3212728
$ws.Cells.Item(110, 10).Value = 318761748
$ws.Cells.Item(110, 11).Value = "کيانمهر_بلواراميرکبير_جنب ميدان امام خميني (ره )"

# Row 111
$ws.Cells.Item(111, 1).Value = "البرز"
$ws.Cells.Item(111, 2).Value = "کرج ناحيه 4"
$ws.Cells.Item(111, 3).Value = "صداقت"
$ws.Cells.Item(111, 4).Value = "دوره ابتدايي توصيفي"
$ws.Cells.Item(111, 5).Value = "پسرانه"
$ws.Cells.Item(111, 6).Value = "غيردولتي"
$ws.Cells.Item(111, 7).Value = "عادي"
$ws.Cells.Item(111, 8).Value = "عادي"
$ws.Cells.Item(111, 9).Value = 34801130
$ws.Cells.Item(111, 10).Value = "-"
$ws.Cells.Item(111, 11).Value = "کرج-خرمدشت-ميثم يک-بهار اول-پلاک37"

# Row 112
$ws.Cells.Item(112, 1).Value = "البرز"
$ws.Cells.Item(112, 2).Value = "کرج ناحيه 4"
$ws.Cells.Item(112, 3).Value = "فرازين"
$ws.Cells.Item(112, 4).Value = "دوره ابتدايي توصيفي"
$ws.Cells.Item(112, 5).Value = "دخترانه"
$ws.Cells.Item(112, 6).Value = "غيردولتي"
$ws.Cells.Item(112, 7).Value = "عادي"
$ws.Cells.Item(112, 8).Value = "عادي"
$ws.Cells.Item(112, 9).Value = 33420115
$ws.Cells.Item(112, 10).Value = "-"
$ws.Cells.Item(112, 11).Value = "کرج-فاز2مهرشهر-بلوار شهرداري-خيابان202-پلاک382/1 پيش و ابتدايي فرازين"

# Row 113
$ws.Cells.Item(113, 1).Value = "البرز"
$ws.Cells.Item(113, 2).Value = "کرج ناحيه 4"
$ws.Cells.Item(113, 3).Value = "حضرت امير(ع)"
$ws.Cells.Item(113, 4).Value = "دوره متوسطه اول"
$ws.Cells.Item(113, 5).Value = "پسرانه"
$ws.Cells.Item(113, 6).Value = "غيردولتي"
$ws.Cells.Item(113, 7).Value = "عادي"
$ws.Cells.Item(113, 8).Value = "عادي"
$ws.Cells.Item(113, 9).Value = 3509123
$ws.Cells.Item(113, 10).Value = 313965359
$ws.Cells.Item(113, 11).Value = "کرج-خيابان درختي-نرسيده به سه راه تهران-پلاک325-متوسطه دوره اول حضرت امير(ع)"

# Row 114
$ws.Cells.Item(114, 1).Value = "البرز"
$ws.Cells.Item(114, 2).Value = "کرج ناحيه 4"
$ws.Cells.Item(114, 3).Value = "شادان"
$ws.Cells.Item(114, 4).Value = "دوره متوسطه اول"
$ws.Cells.Item(114, 5).Value = "دخترانه"
$ws.Cells.Item(114, 6).Value = "غيردولتي"
$ws.Cells.Item(114, 7).Value = "عادي"
$ws.Cells.Item(114, 8).Value = "عادي"
$ws.Cells.Item(114, 9).Value = 34613059
$ws.Cells.Item(114, 10).Value = "-"
$ws.Cells.Item(114, 11).Value = "کرج-حصارک - خيابان برزنت-90دستگاه اول-پلاک34-متوسطه اول شادان"

# Row 115
$ws.Cells.Item(115, 1).Value = "البرز"
$ws.Cells.Item(115, 2).Value = "کرج ناحيه 4"
$ws.Cells.Item(115, 3).Value = "عصر تلاش"
$ws.Cells.Item(115, 4).Value = "دوره ابتدايي توصيفي"
$ws.Cells.Item(115, 5).Value = "پسرانه"
$ws.Cells.Item(115, 6).Value = "غيردولتي"
$ws.Cells.Item(115, 7).Value = "عادي"
$ws.Cells.Item(115, 8).Value = "عادي"
$ws.Cells.Item(115, 9).Value = 34516780
$ws.Cells.Item(115, 10).Value = "-"
$ws.Cells.Item(115, 11).Value = "-"

# Row 116
$ws.Cells.Item(116, 1).Value = "البرز"
$ws.Cells.Item(116, 2).Value = "کرج ناحيه 4"
$ws.Cells.Item(116, 3).Value = "نيوشا"
$ws.Cells.Item(116, 4).Value = "دوره متوسطه اول"
$ws.Cells.Item(116, 5).Value = "دخترانه"
$ws.Cells.Item(116, 6).Value = "غيردولتي"
$ws.Cells.Item(116, 7).Value = "عادي"
$ws.Cells.Item(116, 8).Value = "عادي"
$ws.Cells.Item(116, 9).Value = 33513094
$ws.Cells.Item(116, 10).Value = "-"
$ws.Cells.Item(116, 11).Value = "کرج-خيايان45متري گلشهر-کوچه مينا-پلاک35-متوسطه اول نيوشا"

# Row 117
$ws.Cells.Item(117, 1).Value = "البرز"
$ws.Cells.Item(117, 2).Value = "کرج ناحيه 4"
$ws.Cells.Item(117, 3).Value = "انديشه"
$ws.Cells.Item(117, 4).Value = "متوسطه دوم - هنرستان کاردانش"
$ws.Cells.Item(117, 5).Value = "پسرانه"
$ws.Cells.Item(117, 6).Value = "غيردولتي"
$ws.Cells.Item(117, 7).Value = "عادي"
$ws.Cells.Item(117, 8).Value = "عادي"
$ws.Cells.Item(117, 9).Value = 4641296
$ws.Cells.Item(117, 10).Value = 313891491
$ws.Cells.Item(117, 11).Value = "کرج-45متري گلشهر-آذرشرقي-پلاک14-کاردانش انديشه"

# --- Cosmetic: move the active selection the way the saved workbook shows it ---
$ws.Range("J35").Select()

